$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Summary" header row (old row 39); this shifts old rows 40-45 up to 39-44
$ws.Rows.Item(39).Delete()

# Update column A labels to the new wording/structure
$ws.Range("A1").Value = 'Labels'
$ws.Range("A2").Value = 'Congress'
$ws.Range("A3").Value = 'Session'
$ws.Range("A4").Value = 'Start Date'
$ws.Range("A5").Value = 'End Date'
$ws.Range("A6").Value = 'Civilian '
$ws.Range("A7").Value = '     Civilian, New nominations'
$ws.Range("A8").Value = '     Civilian, Carryover nominations'
$ws.Range("A9").Value = '     Civilian, Confirmed '
$ws.Range("A10").Value = '     Civilian, Withdrawn '
$ws.Range("A11").Value = '     Civilian, Returned to White House '
$ws.Range("A12").Value = 'Other Civilian '
$ws.Range("A13").Value = '     Other Civilian, New nominations'
$ws.Range("A14").Value = '     Other Civilian, Carryover nominations'
$ws.Range("A15").Value = '     Other Civilian, Confirmed '
$ws.Range("A16").Value = '     Other Civilian, Withdrawn '
$ws.Range("A17").Value = '     Other Civilian, Returned to White House '
$ws.Range("A18").Value = 'Air Force '
$ws.Range("A19").Value = '     Air Force, New nominations'
$ws.Range("A20").Value = '     Air Force, Carryover nominations'
$ws.Range("A21").Value = '     Air Force, Confirmed '
$ws.Range("A22").Value = '     Air Force, Withdrawn '
$ws.Range("A23").Value = '     Air Force, Returned to White House '
$ws.Range("A24").Value = 'Army '
$ws.Range("A25").Value = '     Army, New nominations'
$ws.Range("A26").Value = '     Army, Carryover nominations'
$ws.Range("A27").Value = '     Army, Confirmed '
$ws.Range("A28").Value = '     Army, Returned to White House '
$ws.Range("A29").Value = 'Navy '
$ws.Range("A30").Value = '     Navy, New nominations'
$ws.Range("A31").Value = '     Navy, Carryover nominations'
$ws.Range("A32").Value = '     Navy, Confirmed '
$ws.Range("A33").Value = '     Navy, Returned to White House '
$ws.Range("A34").Value = 'Marine Corps '
$ws.Range("A35").Value = '     Marine Corps, New nominations'
$ws.Range("A36").Value = '     Marine Corps, Carryover nominations'
$ws.Range("A37").Value = '     Marine Corps, Confirmed '
$ws.Range("A38").Value = '     Marine Corps, Returned to White House '
$ws.Range("A39").Value = 'Total new nominations'
$ws.Range("A40").Value = 'Total carryover nominations'
$ws.Range("A41").Value = 'Total confirmed '
$ws.Range("A42").Value = 'Total unconfirmed '
$ws.Range("A43").Value = 'Total withdrawn '
$ws.Range("A44").Value = 'Total returned'

# Fix values for rows 39/40 whose figures swap position now that "Summary" is gone
$ws.Range("B39").Value = 24420
$ws.Range("B40").Value = 6812
